$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value2 = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "45.337.02"
Set-TextValue $ws.Range("E2") "  +5.18%  "
Set-TextValue $ws.Range("D3") "2.453.11"
Set-TextValue $ws.Range("E3") "  +3.53%  "
Set-TextValue $ws.Range("E4") "  +0.00%  "
Set-TextValue $ws.Range("D5") "319.92"
Set-TextValue $ws.Range("E5") "  +5.39%  "
Set-TextValue $ws.Range("D6") "104.36"
Set-TextValue $ws.Range("E6") "  +8.34%  "
Set-TextValue $ws.Range("D7") "0.517"
Set-TextValue $ws.Range("E7") "  +2.76%  "
Set-TextValue $ws.Range("E8") "  -0.01%  "
Set-TextValue $ws.Range("D9") "0.533"
Set-TextValue $ws.Range("E9") "  +10.46%  "
Set-TextValue $ws.Range("D10") "36.01"
Set-TextValue $ws.Range("E10") "  +4.67%  "
Set-TextValue $ws.Range("D11") "0.0805"
Set-TextValue $ws.Range("E12") "  -2.16%  "
Set-TextValue $ws.Range("D13") "18.57"
Set-TextValue $ws.Range("E13") "  +1.25%  "
Set-TextValue $ws.Range("E14") "  +3.88%  "
Set-TextValue $ws.Range("D15") "2.839.37"
Set-TextValue $ws.Range("E15") "  +3.77%  "
Set-TextValue $ws.Range("D16") "2.456.78"
Set-TextValue $ws.Range("E16") "  +2.28%  "
Set-TextValue $ws.Range("D17") "0.840"
Set-TextValue $ws.Range("E17") "  +4.75%  "
Set-TextValue $ws.Range("D18") "45.236.25"
Set-TextValue $ws.Range("E18") "  +4.88%  "
Set-TextValue $ws.Range("D19") "12.38"
Set-TextValue $ws.Range("E19") "  +3.24%  "
Set-TextValue $ws.Range("D20") "6.39"
Set-TextValue $ws.Range("E20") "  +1.30%  "
Set-TextValue $ws.Range("E21") "  +4.63%  "
Set-TextValue $ws.Range("D22") "69.24"
Set-TextValue $ws.Range("E22") "  +1.78%  "
Set-TextValue $ws.Range("D23") "244.30"
Set-TextValue $ws.Range("E23") "  +3.64%  "
Set-TextValue $ws.Range("D24") "2.29"
Set-TextValue $ws.Range("E24") "  +2.70%  "
Set-TextValue $ws.Range("D25") "2.53"
Set-TextValue $ws.Range("E25") "  +3.51%  "
Set-TextValue $ws.Range("E26") "  +0.04%  "
Set-TextValue $ws.Range("D27") "25.53"
Set-TextValue $ws.Range("E27") "  +4.34%  "
Set-TextValue $ws.Range("D28") "2.20"
Set-TextValue $ws.Range("E28") "  -7.24%  "
Set-TextValue $ws.Range("E29") "  +2.54%  "
Set-TextValue $ws.Range("D30") "33.91"
Set-TextValue $ws.Range("E30") "  +6.30%  "
Set-TextValue $ws.Range("D31") "49.62"
Set-TextValue $ws.Range("E31") "  +3.42%  "
Set-TextValue $ws.Range("E32") "  +15.04%  "
Set-TextValue $ws.Range("D33") "20.39"
Set-TextValue $ws.Range("E33") "  +14.35%  "
Set-TextValue $ws.Range("E34") "  +4.03%  "
Set-TextValue $ws.Range("E35") "  +0.30%  "
Set-TextValue $ws.Range("D36") "0.0766"
Set-TextValue $ws.Range("E36") "  +3.78%  "
Set-TextValue $ws.Range("E37") "  +5.22%  "
Set-TextValue $ws.Range("E38") "  +4.45%  "
Set-TextValue $ws.Range("D39") "2.88"
Set-TextValue $ws.Range("E39") "  +0.88%  "
Set-TextValue $ws.Range("D40") "125.24"
Set-TextValue $ws.Range("E40") "  -2.37%  "
Set-TextValue $ws.Range("E41") "  +2.49%  "
Set-TextValue $ws.Range("E42") "  -2.59%  "
Set-TextValue $ws.Range("D43") "21.44"
Set-TextValue $ws.Range("E43") "  +1.52%  "
Set-TextValue $ws.Range("D44") "0.0291"
Set-TextValue $ws.Range("E44") "  +4.71%  "
Set-TextValue $ws.Range("D45") "1.945.64"
Set-TextValue $ws.Range("E45") "  +0.86%  "
Set-TextValue $ws.Range("E46") "  +8.59%  "
Set-TextValue $ws.Range("D47") "2.11"
Set-TextValue $ws.Range("E47") "  -1.03%  "
Set-TextValue $ws.Range("E48") "  +0.01%  "
Set-TextValue $ws.Range("D49") "1.78"
Set-TextValue $ws.Range("E49") "  +16.58%  "
Set-TextValue $ws.Range("D50") "76.38"
Set-TextValue $ws.Range("E50") "  +6.51%  "
Set-TextValue $ws.Range("D51") "53.90"
Set-TextValue $ws.Range("E51") "  +4.15%  "
